# Auto-generated edit script updating crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).Value = "57.076.75"
$ws.Cells.Item(2,5).Value = "  +1.23%  "

# Row 3
$ws.Cells.Item(3,4).Value = "2.398.42"
$ws.Cells.Item(3,5).Value = "  +1.57%  "

# Row 4
$ws.Cells.Item(4,5).Value = "  -0.02%  "

# Row 5
$ws.Cells.Item(5,4).Value = "'505.26"
$ws.Cells.Item(5,5).Value = "  -1.44%  "

# Row 6
$ws.Cells.Item(6,4).Value = "'132.95"
$ws.Cells.Item(6,5).Value = "  +4.32%  "

# Row 7
$ws.Cells.Item(7,5).Value = "  +0.01%  "

# Row 8
$ws.Cells.Item(8,5).Value = "  +0.13%  "

# Row 9
$ws.Cells.Item(9,4).Value = "2.410.81"
$ws.Cells.Item(9,5).Value = "  +1.37%  "

# Row 10
$ws.Cells.Item(10,4).Value = "'0.0971"
$ws.Cells.Item(10,5).Value = "  +1.23%  "

# Row 12
$ws.Cells.Item(12,4).Value = "'0.322"
$ws.Cells.Item(12,5).Value = "  +2.02%  "

# Row 13
$ws.Cells.Item(13,4).Value = "'4.57"
$ws.Cells.Item(13,5).Value = "  -4.91%  "

# Row 14
$ws.Cells.Item(14,4).Value = "2.826.02"
$ws.Cells.Item(14,5).Value = "  +1.60%  "

# Row 15
$ws.Cells.Item(15,4).Value = "56.957.70"
$ws.Cells.Item(15,5).Value = "  +1.05%  "

# Row 16
$ws.Cells.Item(16,4).Value = "'21.91"
$ws.Cells.Item(16,5).Value = "  +2.39%  "

# Row 17
$ws.Cells.Item(17,5).Value = "  +2.61%  "

# Row 18
$ws.Cells.Item(18,4).Value = "2.399.03"
$ws.Cells.Item(18,5).Value = "  -0.45%  "

# Row 19
$ws.Cells.Item(19,4).Value = "'10.23"
$ws.Cells.Item(19,5).Value = "  -0.42%  "

# Row 20
$ws.Cells.Item(20,2).Value = "Polkadot"
$ws.Cells.Item(20,3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(20,4).Value = "'4.04"
$ws.Cells.Item(20,5).Value = "  -0.34%  "

# Row 21
$ws.Cells.Item(21,2).Value = "BitcoinCash"
$ws.Cells.Item(21,3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(21,4).Value = "'309.82"
$ws.Cells.Item(21,5).Value = "  -0.07%  "

# Row 22
$ws.Cells.Item(22,4).Value = "'6.33"
$ws.Cells.Item(22,5).Value = "  +2.57%  "

# Row 23
$ws.Cells.Item(23,5).Value = "  +0.54%  "

# Row 24
$ws.Cells.Item(24,4).Value = "'0.997"
$ws.Cells.Item(24,5).Value = "  +0.06%  "

# Row 25
$ws.Cells.Item(25,4).Value = "'64.96"
$ws.Cells.Item(25,5).Value = "  -0.82%  "

# Row 26
$ws.Cells.Item(26,4).Value = "'0.995"
$ws.Cells.Item(26,5).Value = "  -0.55%  "

# Row 27
$ws.Cells.Item(27,2).Value = "Polygon"
$ws.Cells.Item(27,3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(27,4).Value = "'0.376"
$ws.Cells.Item(27,5).Value = "  -3.23%  "

# Row 28
$ws.Cells.Item(28,2).Value = "Kaspa"
$ws.Cells.Item(28,3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(28,4).Value = "'0.153"
$ws.Cells.Item(28,5).Value = "  -0.20%  "

# Row 29
$ws.Cells.Item(29,4).Value = "'7.43"
$ws.Cells.Item(29,5).Value = "  +2.99%  "

# Row 30
$ws.Cells.Item(30,4).Value = "'173.22"
$ws.Cells.Item(30,5).Value = "  -0.92%  "

# Row 31
$ws.Cells.Item(31,4).Value = "0.0₃0724"
$ws.Cells.Item(31,5).Value = "  +1.42%  "

# Row 32
$ws.Cells.Item(32,4).Value = "'1.68"
$ws.Cells.Item(32,5).Value = "  -0.05%  "

# Row 33
$ws.Cells.Item(33,4).Value = "'5.95"
$ws.Cells.Item(33,5).Value = "  -3.06%  "

# Row 34
$ws.Cells.Item(34,5).Value = "  +0.51%  "

# Row 35
$ws.Cells.Item(35,5).Value = "  +0.12%  "

# Row 36
$ws.Cells.Item(36,5).Value = "  +0.03%  "

# Row 37
$ws.Cells.Item(37,5).Value = "  +1.65%  "

# Row 38
$ws.Cells.Item(38,4).Value = "'1.20"
$ws.Cells.Item(38,5).Value = "  +0.71%  "

# Row 39
$ws.Cells.Item(39,5).Value = "  +3.47%  "

# Row 40
$ws.Cells.Item(40,4).Value = "'36.74"

# Row 41
$ws.Cells.Item(41,4).Value = "'0.810"
$ws.Cells.Item(41,5).Value = "  +0.86%  "

# Row 42
$ws.Cells.Item(42,5).Value = "  +1.57%  "

# Row 43
$ws.Cells.Item(43,4).Value = "'132.84"
$ws.Cells.Item(43,5).Value = "  +7.88%  "

# Row 44
$ws.Cells.Item(44,4).Value = "'5.02"
$ws.Cells.Item(44,5).Value = "  +3.26%  "

# Row 45
$ws.Cells.Item(45,4).Value = "'3.36"
$ws.Cells.Item(45,5).Value = "  -0.17%  "

# Row 46
$ws.Cells.Item(46,2).Value = "Bittensor"
$ws.Cells.Item(46,3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(46,4).Value = "'252.78"
$ws.Cells.Item(46,5).Value = "  +0.20%  "

# Row 47
$ws.Cells.Item(47,2).Value = "Mantle"
$ws.Cells.Item(47,3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(47,4).Value = "'0.567"
$ws.Cells.Item(47,5).Value = "  -0.14%  "

# Row 48
$ws.Cells.Item(48,5).Value = "  +0.72%  "

# Row 49
$ws.Cells.Item(49,4).Value = "'0.0489"
$ws.Cells.Item(49,5).Value = "  +0.20%  "

# Row 50
$ws.Cells.Item(50,2).Value = "VeChain"
$ws.Cells.Item(50,3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(50,4).Value = "'0.0210"
$ws.Cells.Item(50,5).Value = "  +1.39%  "

# Row 51
$ws.Cells.Item(51,2).Value = "InjectiveProtocol"
$ws.Cells.Item(51,3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(51,4).Value = "'17.22"
$ws.Cells.Item(51,5).Value = "  +10.21%  "
